# Auto-generated script to apply the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (so Excel does not
# auto-convert numeric-looking strings like "0.999" or "5.70" into numbers),
# then restore the cell to its original (default) style so no stray
# number-format / quote-prefix style is left behind on the cell.
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '68.895.40'
$ws.Range("E2").Value = '  +2.19%  '
$ws.Range("D3").Value = '3.267.66'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue "D5" '584.47'
$ws.Range("E5").Value = '  +1.14%  '
Set-TextValue "D6" '181.81'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  -0.08%  '
Set-TextValue "D8" '0.597'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +1.51%  '
$ws.Range("E10").Value = '  -0.90%  '
Set-TextValue "D11" '0.423'
$ws.Range("E11").Value = '  +2.03%  '
$ws.Range("D12").Value = '3.835.18'
$ws.Range("E12").Value = '  +0.56%  '
Set-TextValue "D14" '28.46'
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").Value = '68.801.20'
$ws.Range("E15").Value = '  +2.07%  '
Set-TextValue "D16" '0.0000171'
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("D17").Value = '3.219.65'
$ws.Range("E17").Value = '  -1.08%  '
Set-TextValue "D18" '5.83'
$ws.Range("E18").Value = '  -0.30%  '
Set-TextValue "D19" '13.53'
$ws.Range("E19").Value = '  -0.25%  '
Set-TextValue "D20" '393.26'
$ws.Range("E20").Value = '  +4.80%  '
Set-TextValue "D21" '7.67'
$ws.Range("E21").Value = '  +0.69%  '
Set-TextValue "D22" '71.83'
$ws.Range("E22").Value = '  +1.02%  '
Set-TextValue "D23" '0.999'
$ws.Range("E23").Value = '  -0.10%  '
Set-TextValue "D24" '0.516'
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("E26").Value = '  +5.77%  '
Set-TextValue "D27" '9.59'
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  +0.85%  '
Set-TextValue "D30" '5.70'
$ws.Range("E30").Value = '  -0.55%  '
Set-TextValue "D31" '22.98'
$ws.Range("E31").Value = '  +1.32%  '
Set-TextValue "D32" '7.12'
$ws.Range("E32").Value = '  +3.01%  '
Set-TextValue "D33" '1.28'
$ws.Range("E33").Value = '  +0.95%  '
Set-TextValue "D35" '163.87'
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("E36").Value = '  -0.09%  '
Set-TextValue "D37" '1.94'
$ws.Range("E37").Value = '  +5.14%  '
Set-TextValue "D38" '0.828'
$ws.Range("E38").Value = '  -2.49%  '
Set-TextValue "D39" '26.36'
$ws.Range("E39").Value = '  -1.05%  '
Set-TextValue "D40" '4.59'
$ws.Range("E40").Value = '  -1.69%  '
Set-TextValue "D41" '6.61'
$ws.Range("E41").Value = '  -3.37%  '
Set-TextValue "D42" '41.35'
$ws.Range("E42").Value = '  +1.28%  '
Set-TextValue "D43" '2.48'
$ws.Range("E43").Value = '  -5.16%  '
Set-TextValue "D44" '0.0689'
$ws.Range("E44").Value = '  +1.30%  '
Set-TextValue "D45" '346.56'
$ws.Range("E45").Value = '  -3.14%  '
$ws.Range("D46").Value = '2.605.09'
$ws.Range("E46").Value = '  -3.72%  '
Set-TextValue "D47" '24.71'
$ws.Range("E47").Value = '  -2.79%  '
$ws.Range("E48").Value = '  +0.95%  '
$ws.Range("E49").Value = '  +2.59%  '
Set-TextValue "D50" '31.61'
$ws.Range("E50").Value = '  +1.49%  '
$ws.Range("E51").Value = '  -0.51%  '
